$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Import statement in cell B2 (merged B2:F2) to reference the
# new package for RequirementForm / CaseData.
$ws.Range("B2").Value = "com.redhat.demos.dm.loan.model.RequirementForm,`ncom.redhat.demos.dm.loan.model.CaseData"

# Reflect the new selection (the edited cell / merged range) as recorded
# in the saved workbook view.
$ws.Range("B2:F2").Select()
